# Updated slides and code examples for synchronous operations
#
# The deck's "date" placeholder (the auto-updating datetimeFigureOut field
# that lives on the slide master and on every slide layout) was showing a
# stale cached value ("23/11/2018"). Refresh it to "08/03/2020" everywhere
# it appears - on the slide master itself and on each of its custom
# (slide) layouts.

$p = $ppt.ActivePresentation
$newDate = "08/03/2020"

function Update-DatePlaceholders {
    param($shapes, [string]$text)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.Name -like "Date Placeholder*") {
            $shape.TextFrame.TextRange.Text = $text
        }
    }
}

# Slide master's own date placeholder.
Update-DatePlaceholders $p.SlideMaster.Shapes $newDate

# Every slide layout hanging off the master has its own date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholders $layout.Shapes $newDate
}
